# B1--and-B2-PowerPoint.pptx edit
#
# 1) Slide 5 contains a table ("Google Shape;122;p17", Shapes.Item(2)) whose
#    table style is changed from the deck's custom style
#    {94A321BE-5E99-40B3-8094-F1EE990CCAA3} to the built-in style
#    {E37124FB-D26E-4B1E-911C-D039011C1263}.
#
# 2) The presentation's active theme (ppt/theme/theme2.xml, the theme used
#    by the slide master / all slides) is switched from the "Integral" /
#    "Red Violet" palette to the "Office Theme" / "Office" palette - i.e.
#    the same 12 theme colors that used to live in ppt/theme/theme2.xml are
#    replaced with the colors that used to live in ppt/theme/theme1.xml
#    (font scheme / format scheme are identical between the two themes, so
#    only the color scheme actually changes visually).

$p = $ppt.ActivePresentation

# --- 1) Table style -------------------------------------------------------
$slide = $p.Slides.Item(5)
$tableShape = $slide.Shapes.Item(2)
$table = $tableShape.Table
$table.ApplyStyle("{E37124FB-D26E-4B1E-911C-D039011C1263}")

# --- 2) Theme colors: Integral/Red Violet -> Office Theme/Office ---------
# ThemeColorScheme.Colors index order: dk1, lt1, dk2, lt2,
# accent1..accent6, hlink, folHlink. RGB values below are the "Office"
# scheme's colors, encoded as COM BGR-packed integers (R | G<<8 | B<<16).
$slide1 = $p.Slides.Item(1)
$themeColors = $slide1.ThemeColorScheme

$officeRGB = @(
    0,          # dk1     000000
    16777215,   # lt1     FFFFFF
    6968388,    # dk2     44546A
    15132391,   # lt2     E7E6E6
    13998939,   # accent1 5B9BD5
    3243501,    # accent2 ED7D31
    10855845,   # accent3 A5A5A5
    49407,      # accent4 FFC000
    12874308,   # accent5 4472C4
    4697456,    # accent6 70AD47
    12673797,   # hlink   0563C1
    7491477     # folHlink 954F72
)

for ($i = 1; $i -le $themeColors.Count; $i++) {
    $themeColors.Colors($i).RGB = $officeRGB[$i - 1]
}
